$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 72, shifting existing rows 72:156 down to 73:157
$ws.Rows.Item(72).Insert()

# Fill in the new row 72 with the new record's data
$ws.Range("A72").Value = 10
$ws.Range("B72").Value = "Vega Modelo de Temuco"
$ws.Range("C72").Value = "La Araucanía"
$ws.Range("D72").Value = 44671
$ws.Range("E72").Value = 9
$ws.Range("F72").Value = 100114007
$ws.Range("G72").Value = "Jengibre"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 15
$ws.Range("K72").Value = 20000
$ws.Range("L72").Value = 20000
$ws.Range("M72").Value = 20000
$ws.Range("N72").Value = "$/caja 13 kilos"
$ws.Range("O72").Value = "Perú"
$ws.Range("P72").Value = 1538
$ws.Range("Q72").Value = 13
$ws.Range("R72").Value = "Hortaliza"

# Match the date-format style used by the rest of column D
$ws.Range("D72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
